# Update the "Förändrad" (Changed) date column (C) from 2023-10-04 (45203)
# to 2023-10-05 (45204) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the sheet's used range.
$lastRow = $ws.UsedRange.Rows.Count

# Data starts on row 2 (row 1 holds the header "Förändrad").
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45204
    }
}
